# Apply the "trips" export template changes:
#  - Pass client timezone into the date/time formulas
#  - Switch OpenStreetMap links from http to https
#  - Update the active selection to D9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Period: value (row 6, column B) - now uses from/to .toString(...) instead of "".format(...)
$ws.Range("B6").Value2 = '${from.toString("YYYY.MM.dd HH:mm:ss")+" - "+to.toString("YYYY.MM.dd HH:mm:ss")}'

# Trip start time (row 9, column A) - now wraps trip.startTime in a joda DateTime with timezone
$ws.Range("A9").Value2 = '${new("org.joda.time.DateTime", trip.startTime, timezone).toString("YYYY.MM.dd HH:mm:ss")}'

# Trip end time (row 9, column C) - now wraps trip.endTime in a joda DateTime with timezone
$ws.Range("C9").Value2 = '${new("org.joda.time.DateTime", trip.endTime, timezone).toString("YYYY.MM.dd HH:mm:ss")}'

# Trip start address hyperlink (row 9, column B) - http -> https
$ws.Range("B9").Value2 = '${util.hyperlink("".format("https://www.openstreetmap.org/?mlat=%1$f&mlon=%2$f#map=16/%1$f/%2$f", trip.startLat, trip.startLon), trip.getStartAddress() == null ? "".format("%1$f°, %2$f°", trip.startLat, trip.startLon) : trip.startAddress)}'

# Trip end address hyperlink (row 9, column D) - http -> https
$ws.Range("D9").Value2 = '${util.hyperlink("".format("https://www.openstreetmap.org/?mlat=%1$f&mlon=%2$f#map=16/%1$f/%2$f", trip.endLat, trip.endLon), trip.getEndAddress() == null ? "".format("%1$f°, %2$f°", trip.endLat, trip.endLon) : trip.endAddress)}'

# Move the active selection from B2 to D9
$ws.Range("D9").Select() | Out-Null
